$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-24 Tuesday" "2025-06-25 Wednesday"
Replace-Text "294×5=" "999×9="
Replace-Text "473×6=" "906×3="
Replace-Text "401×3=" "405×9="
Replace-Text "344×5=" "743×2="
Replace-Text "808×7=" "724×4="
Replace-Text "572×6=" "254×4="
Replace-Text "686×6=" "492×8="
Replace-Text "336×9=" "810×8="
Replace-Text "424×6=" "159×8="
Replace-Text "342×4=" "317×6="
Replace-Text "549×2=" "558×2="
Replace-Text "748×7=" "947×4="
Replace-Text "414×4=" "998×4="
Replace-Text "158×6=" "373×2="
Replace-Text "246×4=" "184×8="
Replace-Text "815×8=" "474×7="
Replace-Text "553×4=" "762×4="
Replace-Text "800×4=" "679×3="
Replace-Text "789×5=" "488×9="
Replace-Text "532×6=" "568×8="
Replace-Text "852×6=" "499×8="
Replace-Text "949×2=" "352×7="
Replace-Text "364×8=" "793×7="
Replace-Text "408×9=" "787×2="
Replace-Text "998×7=" "457×5="
